$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header changes
$ws.Range("K1").Value = "Tanggal"
$ws.Range("M1").Value = "Skema sidang"

# Password column (O) regeneration for rows 2-274
$passwords = @{
    2 = "3kQhCo"
    3 = "qmCLTU"
    4 = "HSGD4N"
    5 = "o2pwOk"
    6 = "ijpzyJ"
    7 = "lurvbj"
    8 = "LPu7hX"
    9 = "rDnbD3"
    10 = "3Prn8e"
    11 = "jJXkbU"
    12 = "wW8v9E"
    13 = "wE8Ubb"
    14 = "28TBw1"
    15 = "ggMqtz"
    16 = "RRY2TP"
    17 = "mrXpTo"
    18 = "OcxmId"
    19 = "KOHkNO"
    20 = "AReLo5"
    21 = "H7wN7H"
    22 = "az8BZp"
    23 = "vlZOYN"
    24 = "wjOkb2"
    25 = "AiSjJO"
    26 = "LEFjiw"
    27 = "vqYdDM"
    28 = "KkZZzG"
    29 = "VuKEAS"
    30 = "qRyat0"
    31 = "qkAZZZ"
    32 = "aZ0WPs"
    33 = "pe5jBe"
    34 = "1RX0Rq"
    35 = "XOrExp"
    36 = "AB2UwC"
    37 = "REJkN6"
    38 = "3cAY8A"
    39 = "s8Lf1g"
    40 = "jwMXi1"
    41 = "AXEy9s"
    42 = "hzwgUT"
    43 = "DRWo1A"
    44 = "WWabHH"
    45 = "8wADs8"
    46 = "shCAsR"
    47 = "zraMCK"
    48 = "Zg49xY"
    49 = "vyUmJN"
    50 = "5Yv17H"
    51 = "GqaXex"
    52 = "dcVBWU"
    53 = "hFxvuW"
    54 = "JdboPx"
    55 = "99pKbd"
    56 = "nAh4I2"
    57 = "c0UaGr"
    58 = "iQecZ8"
    59 = "5ENQm2"
    60 = "wejnNa"
    61 = "HObtuE"
    62 = "RUoX3b"
    63 = "Bhj5ox"
    64 = "QLDVTe"
    65 = "nxYP4F"
    66 = "reIkbp"
    67 = "XNn0BI"
    68 = "G6uoDd"
    69 = "KEkZTa"
    70 = "BSUuAV"
    71 = "lWEnHH"
    72 = "V04sHj"
    73 = "zwSFZ1"
    74 = "61AFgG"
    75 = "wUmbYH"
    76 = "X9h29Z"
    77 = "bUojnp"
    78 = "6rv3oP"
    79 = "4iMlIU"
    80 = "9W6MsP"
    81 = "ueQow7"
    82 = "EvNTEu"
    83 = "06aZuz"
    84 = "Fnb9AE"
    85 = "hUoy3C"
    86 = "B7ecFR"
    87 = "tXCiqO"
    88 = "SCFB3H"
    89 = "OF7ce1"
    90 = "WttUhh"
    91 = "FpgKto"
    92 = "tZHpWW"
    93 = "KuBGFL"
    94 = "wy0S69"
    95 = "TECa4Q"
    96 = "XKkmu5"
    97 = "UHE1Aa"
    98 = "cQT2LC"
    99 = "8c5r5c"
    100 = "kZOyzf"
    101 = "JndCxH"
    102 = "tpHZRl"
    103 = "dQ3AdY"
    104 = "t9rBLi"
    105 = "06XvGU"
    106 = "0riTVW"
    107 = "Y3kv2P"
    108 = "82kNNL"
    109 = "VzqaMT"
    110 = "InNQez"
    111 = "XeAdsC"
    112 = "NOyJya"
    113 = "RIoQVH"
    114 = "nwg44n"
    115 = "VRJyqQ"
    116 = "6exf1Q"
    117 = "sTChug"
    118 = "mGYO2P"
    119 = "rHM42w"
    120 = "8vTQtl"
    121 = "XKjtVj"
    122 = "aVjVEs"
    123 = "EB27Wb"
    124 = "tb0YwW"
    125 = "tEihot"
    126 = "ZvH0eF"
    127 = "MfFjpJ"
    128 = "EzXsmB"
    129 = "E3AqT4"
    130 = "h4b4pp"
    131 = "PqFj1U"
    132 = "bu7cPR"
    133 = "tP9hoz"
    134 = "qJt4wm"
    135 = "ryQybw"
    136 = "c3PnL7"
    137 = "0ooVV0"
    138 = "kVIIOf"
    139 = "3AOCX4"
    140 = "bO9uhg"
    141 = "BIGkdM"
    142 = "eabImJ"
    143 = "8Vf5RA"
    144 = "yWKURD"
    145 = "pYhvZc"
    146 = "DvaKqy"
    147 = "XdIa7L"
    148 = "RkNV3q"
    149 = "MzNj37"
    150 = "m5YYug"
    151 = "K6Xcp6"
    152 = "f0RKHU"
    153 = "lnmfvh"
    154 = "7RXvZI"
    155 = "4z6vW6"
    156 = "TamS55"
    157 = "fRiR4T"
    158 = "cI3RqI"
    159 = "RELgdd"
    160 = "L2a9WZ"
    161 = "JRCzTG"
    162 = "0kA89S"
    163 = "HlhQhX"
    164 = "aAKgfw"
    165 = "BJJz5P"
    166 = "TjEFrv"
    167 = "9PAFXD"
    168 = "bGF8MB"
    169 = "DlAX8K"
    170 = "cCRKLM"
    171 = "DXzYWy"
    172 = "nsjcaa"
    173 = "k5YWPq"
    174 = "JapNQn"
    175 = "qKfwjb"
    176 = "1vXJd0"
    177 = "rDCHjY"
    178 = "UZpO0a"
    179 = "fj6P4G"
    180 = "58UHZY"
    181 = "bxiug9"
    182 = "um6rQP"
    183 = "YzqXnH"
    184 = "GztFAS"
    185 = "OmcL44"
    186 = "UV4Q3D"
    187 = "T6lTmK"
    188 = "DzkX9y"
    189 = "IhtwbO"
    190 = "q5ZvDC"
    191 = "YIzDIo"
    192 = "vs8a9T"
    193 = "cRkUGC"
    194 = "3u8bnF"
    195 = "h9WuXD"
    196 = "KrBo5l"
    197 = "LYCogq"
    198 = "YVOfSG"
    199 = "3evd3v"
    200 = "61nYYL"
    201 = "Y6g1WV"
    202 = "uZ0Lu2"
    203 = "St4JZ5"
    204 = "v3tXSW"
    205 = "8yR4KP"
    206 = "9iva2t"
    207 = "c8iuMT"
    208 = "41VED9"
    209 = "i3JJso"
    210 = "uMXA7s"
    211 = "nxTCUQ"
    212 = "lzAf2y"
    213 = "KuP8Cg"
    214 = "mLza5g"
    215 = "iznNJQ"
    216 = "D7lAaL"
    217 = "tVutMI"
    218 = "wgfvUT"
    219 = "3gVaff"
    220 = "Q7bFhf"
    221 = "nACsVR"
    222 = "loNekX"
    223 = "G1hUzz"
    224 = "yOqH3C"
    225 = "Dw2FYO"
    226 = "huahkK"
    227 = "nOBRL9"
    228 = "x8j9W7"
    229 = "HOix2J"
    230 = "bjP793"
    231 = "eGCIHt"
    232 = "gmblZo"
    233 = "HEl2c6"
    234 = "UXOZrD"
    235 = "xkl2Ah"
    236 = "5HtU0x"
    237 = "N63YpR"
    238 = "Jaidun"
    239 = "chC8gn"
    240 = "ImtphX"
    241 = "k3b2rD"
    242 = "YOlRR9"
    243 = "nQUbg0"
    244 = "dVTvzi"
    245 = "5yBfi9"
    246 = "MSCrqu"
    247 = "G2ydqg"
    248 = "ZXRpUC"
    249 = "ACKRD8"
    250 = "JGp6bi"
    251 = "3owte0"
    252 = "IglZKc"
    253 = "iTa1Sr"
    254 = "iko9lg"
    255 = "2ngGPv"
    256 = "xdrRpR"
    257 = "OrpBuX"
    258 = "Q4RZcD"
    259 = "xvueBf"
    260 = "LmLSYJ"
    261 = "lU4kas"
    262 = "5NzZZR"
    263 = "2hQDQF"
    264 = "HlCbmv"
    265 = "2imFK6"
    266 = "o0vcZk"
    267 = "fQg08A"
    268 = "eFIZpD"
    269 = "z3uHzh"
    270 = "T4V0Ce"
    271 = "AiM1ii"
    272 = "GL27nn"
    273 = "q3F3Ty"
    274 = "LeMVsX"
}

foreach ($row in $passwords.Keys) {
    $ws.Cells.Item($row, 15).Value = $passwords[$row]
}

Write-Output "Done"